$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation is inserted immediately after row 71.
# The existing rows 72/73/74 shift down to become rows 73/74/75, and the
# new observation is written into row 72.

$cols = 1..18   # columns A..R

# Shift rows downward, starting from the bottom so source data for a
# later copy isn't clobbered before it is read.
foreach ($c in $cols) {
    $ws.Cells.Item(75, $c).Value2 = $ws.Cells.Item(74, $c).Value2
}
foreach ($c in $cols) {
    $ws.Cells.Item(74, $c).Value2 = $ws.Cells.Item(73, $c).Value2
}
foreach ($c in $cols) {
    $ws.Cells.Item(73, $c).Value2 = $ws.Cells.Item(72, $c).Value2
}

# Row 75 is a brand new row so it doesn't inherit the date-formatted style
# used by column D (s="2" / numFmt "YYYY-MM-DD HH:MM:SS"); copy that over.
$ws.Cells.Item(75, 4).NumberFormat = $ws.Cells.Item(74, 4).NumberFormat

# Write the new record into row 72.
$ws.Cells.Item(72, 1).Value2 = 10
$ws.Cells.Item(72, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(72, 3).Value2 = "La Araucanía"
$ws.Cells.Item(72, 4).Value2 = 44610
$ws.Cells.Item(72, 5).Value2 = 9
$ws.Cells.Item(72, 6).Value2 = 100112030
$ws.Cells.Item(72, 7).Value2 = "Poroto granado"
$ws.Cells.Item(72, 8).Value2 = "Sin especificar"
$ws.Cells.Item(72, 9).Value2 = "Primera"
$ws.Cells.Item(72, 10).Value2 = 50
$ws.Cells.Item(72, 11).Value2 = 28000
$ws.Cells.Item(72, 12).Value2 = 28000
$ws.Cells.Item(72, 13).Value2 = 28000
$ws.Cells.Item(72, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(72, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(72, 16).Value2 = 1120
$ws.Cells.Item(72, 17).Value2 = 25
$ws.Cells.Item(72, 18).Value2 = "Hortaliza"
